$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    99.430199430199423,
    99.430199430199423,
    99.430199430199423,
    99.470899470899468,
    99.470899470899468,
    99.470899470899468,
    99.26739926739927,
    99.26739926739927,
    99.26739926739927,
    99.26739926739927,
    98.697598697598693,
    99.26739926739927,
    99.348799348799346,
    99.26739926739927,
    99.348799348799346,
    99.26739926739927,
    99.26739926739927,
    99.26739926739927,
    99.674399674399666,
    99.430199430199423,
    99.430199430199423,
    99.430199430199423,
    99.430199430199423,
    99.430199430199423,
    99.308099308099301,
    99.26739926739927,
    99.308099308099301,
    99.348799348799346,
    99.26739926739927,
    99.348799348799346,
    99.430199430199423,
    99.470899470899468,
    99.430199430199423,
    99.470899470899468,
    99.430199430199423,
    99.348799348799346,
    99.470899470899468,
    99.470899470899468,
    99.470899470899468,
    99.26739926739927,
    99.26739926739927,
    99.348799348799346,
    99.348799348799346,
    99.348799348799346,
    99.348799348799346,
    99.348799348799346,
    99.26739926739927,
    99.348799348799346
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
